$wb = $excel.ActiveWorkbook

# --- Sheet "Nädal 5": just a view/selection change (cosmetic) ---
$ws5 = $wb.Worksheets.Item(5)

# --- Sheet "Nädal 6": remove the 5 unused empty rows (10-14) and complete the
#     last filled-in entry (task 22 finished, task 23 started) ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows("16:20").Delete()

# --- Duplicate "Nädal 6" (after the row cleanup above) to create the new
#     week's sheet, "Nädal 7", placed right after it ---
$ws6.Copy([System.Reflection.Missing]::Value, $ws6)
$ws7 = $wb.Worksheets.Item(7)
$ws7.Name = "Nädal 7"

# New week header date range on the freshly duplicated sheet.
$ws7.Range("G4").Value = "10.03.2020 - 16.03.2020"

# Fill in the remaining data for row 15 on "Nädal 6" (task 22 wrap-up).
$ws6.Range("D15").Value = 0.70833333333333337
$ws6.Range("F15").Value = 60
$ws6.Range("H15").Value = "p. 22 tehtud, alustanud p. 23"

# Clear the data that got duplicated onto "Nädal 7" (keep the row index
# numbers in column A, clear everything else) so the new week starts blank.
$ws7.Range("B7:J15").ClearContents()

# Log the first entry of week 7.
$ws7.Range("B7").Value = 43893
$ws7.Range("C7").Value = 0.82986111111111116
$ws7.Range("D7").Value = 0.87777777777777777
$ws7.Range("F7").Value = 69
$ws7.Range("G7").Value = "Kodutöö 6"
$ws7.Range("H7").Value = "p. 23 tehtud"
$ws7.Range("J7").Value = "x"

# --- Restore the view/selection state for each touched sheet ---
$ws5.Activate()
$ws5.Range("G19").Select()

$ws6.Activate()
$ws6.Range("H17").Select()

$ws7.Activate()
$ws7.Range("H7").Select()
